$p = $ppt.ActivePresentation
$m = $p.SlideMaster
$cs = $m.ColorScheme
for ($i=1;$i -le $cs.Count;$i++) {
  Write-Output "SM Item $i : $($cs.Item($i).RGB)"
}
$nm = $p.NotesMaster
$cs2 = $nm.ColorScheme
for ($i=1;$i -le $cs2.Count;$i++) {
  Write-Output "NM Item $i : $($cs2.Item($i).RGB)"
}
